# Refresh the crypto price/volume snapshot on Sheet1 (also fixes the
# Theta/ApeX rank swap between rows 46 and 47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores plain-text numbers (e.g. "607.22", "70.290.96").
# A General-formatted cell receiving a numeric-looking string is auto-converted
# to a real number by Excel, which would both change the cells stored type and,
# for values with a significant trailing zero such as "186.20" or "9.20", drop
# that trailing zero on display. Temporarily marking those specific cells as Text
# (@) before assigning keeps the new price strings verbatim; clearing the format
# again afterwards restores the original (default) cell style.
$textLockCells = @("D5", "D6", "D12", "D14", "D18", "D19", "D20", "D22", "D23", "D26", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D46", "D47", "D49", "D51")
foreach ($addr in $textLockCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.290.96"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "3.561.77"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "607.22"
$ws.Range("E5").Value = "  +3.33%  "

$ws.Range("D6").Value = "186.20"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").Value = "3.552.43"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("E10").Value = "  +8.99%  "

$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "53.88"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").Value = "9.54"
$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").Value = "4.124.84"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").Value = "70.381.91"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "3.579.85"
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "12.75"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("D19").Value = "19.02"
$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("D20").Value = "579.16"
$ws.Range("E20").Value = "  +7.29%  "

$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "0.993"
$ws.Range("E22").Value = "  -2.13%  "

$ws.Range("D23").Value = "17.35"
$ws.Range("E23").Value = "  -3.64%  "

$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("E25").Value = "  +0.32%  "

$ws.Range("D26").Value = "94.28"
$ws.Range("E26").Value = "  -1.33%  "

$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").Value = "10.93"
$ws.Range("E28").Value = "  -2.95%  "

$ws.Range("D29").Value = "9.42"
$ws.Range("E29").Value = "  +3.20%  "

$ws.Range("D30").Value = "32.29"
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("D31").Value = "7.07"
$ws.Range("E31").Value = "  -3.29%  "

$ws.Range("D32").Value = "12.23"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D34").Value = "63.62"
$ws.Range("E34").Value = "  -1.40%  "

$ws.Range("D35").Value = "3.68"
$ws.Range("E35").Value = "  +18.48%  "

$ws.Range("D36").Value = "3.18"
$ws.Range("E36").Value = "  -1.92%  "

$ws.Range("D37").Value = "528.77"
$ws.Range("E37").Value = "  -3.82%  "

$ws.Range("D38").Value = "0.404"
$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").Value = "37.42"
$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("D41").Value = "0.0₃0787"
$ws.Range("E41").Value = "  +2.60%  "

$ws.Range("D42").Value = "3.531.51"
$ws.Range("E42").Value = "  +5.67%  "

$ws.Range("E43").Value = "  +4.33%  "

$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("E45").Value = "  +3.94%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.46"
$ws.Range("E46").Value = "  -3.54%  "

$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("E48").Value = "  +3.45%  "

$ws.Range("D49").Value = "9.20"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("D51").Value = "136.73"
$ws.Range("E51").Value = "  -0.39%  "

foreach ($addr in $textLockCells) {
    $ws.Range($addr).ClearFormats()
}
